$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Order of writes matters for shared-string table ordering; matches
# original authoring sequence.

# Cumulative total label updated (first new string)
$ws.Range("A14").Value = "Cumulative Total: 200"

# Row 3 - Stage
$ws.Range("A3").Value = "Project Build"

# Row 5 - Stage / Task
$ws.Range("A5").Value = "Proj analysis/elicitation"
$ws.Range("B5").Value = "Interview with client, redefine requirements from feedback"

# Row 6 - Task
$ws.Range("B6").Value = "Work on final iteration"

# Header - name
$ws.Range("C1").Value = "Richard Dobson"

# Row 3 - Task
$ws.Range("B3").Value = "Work on sprint"

# Row 4 - Task
$ws.Range("B4").Value = "Complete sprint"

# Remaining cells reusing existing shared strings / plain numbers
$ws.Range("A4").Value = "Project Build"
$ws.Range("A6").Value = "Project Build"

$ws.Range("E1").Value = 10

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 7

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 7

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3

# Column width change
$ws.Columns.Item(1).ColumnWidth = 23.5703125

# Selection change
$ws.Range("B5").Select()

$wb.Save()
